# Applies the "Updated cryptos list" refresh: new coin name/link for rows
# 49-51 (EnergySwap/Maker/VeChain -> Maker/VeChain/FLOKI) and refreshed
# Price (D) / Volume(1h) (E) figures for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    # Force literal text so numeric-looking strings (e.g. "1.00", "68.512.04")
    # are not reinterpreted as numbers, then restore the default "Normal"
    # cell style so no stray formatting is introduced.
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "68.512.04"
$ws.Range("E2").Value = "  -0.96%  "
# Row 3
Set-TextValue $ws.Range("D3") "3.898.34"
$ws.Range("E3").Value = "  +2.58%  "
# Row 4
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.07%  "
# Row 5
Set-TextValue $ws.Range("D5") "602.22"
$ws.Range("E5").Value = "  +0.02%  "
# Row 6
Set-TextValue $ws.Range("D6") "166.05"
$ws.Range("E6").Value = "  +1.25%  "
# Row 7
Set-TextValue $ws.Range("D7") "3.895.43"
$ws.Range("E7").Value = "  +2.59%  "
# Row 8
$ws.Range("E8").Value = "  +0.07%  "
# Row 9
$ws.Range("E9").Value = "  -1.32%  "
# Row 10
$ws.Range("E10").Value = "  -1.68%  "
# Row 11
$ws.Range("E11").Value = "  +1.36%  "
# Row 12
$ws.Range("E12").Value = "  -0.26%  "
# Row 13
Set-TextValue $ws.Range("D13") "0.0000256"
$ws.Range("E13").Value = "  +4.24%  "
# Row 14
Set-TextValue $ws.Range("D14") "37.28"
$ws.Range("E14").Value = "  -0.03%  "
# Row 15
Set-TextValue $ws.Range("D15") "4.551.33"
$ws.Range("E15").Value = "  +2.76%  "
# Row 16
Set-TextValue $ws.Range("D16") "3.912.11"
$ws.Range("E16").Value = "  +3.21%  "
# Row 17
Set-TextValue $ws.Range("D17") "68.594.81"
$ws.Range("E17").Value = "  -0.93%  "
# Row 18
Set-TextValue $ws.Range("D18") "7.48"
$ws.Range("E18").Value = "  +0.81%  "
# Row 19
Set-TextValue $ws.Range("D19") "17.18"
$ws.Range("E19").Value = "  -0.72%  "
# Row 20
$ws.Range("E20").Value = "  -2.35%  "
# Row 21
Set-TextValue $ws.Range("D21") "11.02"
$ws.Range("E21").Value = "  -3.31%  "
# Row 22
Set-TextValue $ws.Range("D22") "486.98"
$ws.Range("E22").Value = "  -0.26%  "
# Row 23
Set-TextValue $ws.Range("D23") "0.723"
$ws.Range("E23").Value = "  -0.01%  "
# Row 24
$ws.Range("E24").Value = "  +10.68%  "
# Row 25
Set-TextValue $ws.Range("D25") "84.53"
$ws.Range("E25").Value = "  -0.19%  "
# Row 26
Set-TextValue $ws.Range("D26") "2.23"
$ws.Range("E26").Value = "  -1.27%  "
# Row 27
$ws.Range("E27").Value = "  -1.58%  "
# Row 28
Set-TextValue $ws.Range("D28") "10.11"
$ws.Range("E28").Value = "  +0.41%  "
# Row 29
$ws.Range("E29").Value = "  +0.07%  "
# Row 30
$ws.Range("E30").Value = "  -0.99%  "
# Row 31
Set-TextValue $ws.Range("D31") "4.050.44"
$ws.Range("E31").Value = "  +2.71%  "
# Row 32
$ws.Range("E32").Value = "  -0.96%  "
# Row 33
$ws.Range("E33").Value = "  -4.25%  "
# Row 34
Set-TextValue $ws.Range("D34") "31.82"
$ws.Range("E34").Value = "  +0.14%  "
# Row 35
Set-TextValue $ws.Range("D35") "3.853.57"
$ws.Range("E35").Value = "  +3.01%  "
# Row 36
$ws.Range("E36").Value = "  -0.35%  "
# Row 37
$ws.Range("E37").Value = "  +2.05%  "
# Row 38
Set-TextValue $ws.Range("D38") "5.93"
$ws.Range("E38").Value = "  +0.11%  "
# Row 39
$ws.Range("E39").Value = "  -1.98%  "
# Row 40
Set-TextValue $ws.Range("D40") "3.20"
$ws.Range("E40").Value = "  +6.29%  "
# Row 41
Set-TextValue $ws.Range("D41") "0.999"
$ws.Range("E41").Value = "  -0.07%  "
# Row 42
$ws.Range("E42").Value = "  -2.32%  "
# Row 43
Set-TextValue $ws.Range("D43") "429.44"
$ws.Range("E43").Value = "  +2.04%  "
# Row 44
Set-TextValue $ws.Range("D44") "1.98"
$ws.Range("E44").Value = "  -0.96%  "
# Row 45
Set-TextValue $ws.Range("D45") "48.30"
$ws.Range("E45").Value = "  -0.62%  "
# Row 46
Set-TextValue $ws.Range("D46") "8.51"
$ws.Range("E46").Value = "  +1.60%  "
# Row 48
Set-TextValue $ws.Range("D48") "142.28"
$ws.Range("E48").Value = "  +0.67%  "
# Row 49
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D49") "2.808.26"
$ws.Range("E49").Value = "  -0.62%  "
# Row 50
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D50") "0.0352"
$ws.Range("E50").Value = "  +0.61%  "
# Row 51
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
Set-TextValue $ws.Range("D51") "0.000264"
$ws.Range("E51").Value = "  +18.00%  "
